$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (precision)
$ws.Range("C2").Value = 0.4
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("I2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3333333333333333
$ws.Range("Q2").Value = 0.4
$ws.Range("S2").Value = 0.5
$ws.Range("W2").Value = 0.5

# Row 3 (recall)
$ws.Range("C3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("S3").Value = 0.5
$ws.Range("W3").Value = 0.5

# Row 4 (f1-score)
$ws.Range("C4").Value = 0.5714285714285715
$ws.Range("F4").Value = 0.5
$ws.Range("I4").Value = 0.8
$ws.Range("M4").Value = 0.5
$ws.Range("Q4").Value = 0.5714285714285715
$ws.Range("S4").Value = 0.5
$ws.Range("W4").Value = 0.5

# Row 5 (f2-score)
$ws.Range("C5").Value = 0.7692307692307692
$ws.Range("F5").Value = 0.7142857142857143
$ws.Range("I5").Value = 0.9090909090909091
$ws.Range("M5").Value = 0.7142857142857143
$ws.Range("Q5").Value = 0.7692307692307692
$ws.Range("S5").Value = 0.5
$ws.Range("W5").Value = 0.5

# Row 6 (NDCG)
$ws.Range("C6").Value = 0.9639404333166532
$ws.Range("F6").Value = 0.944847956559586
$ws.Range("I6").Value = 0.7967075809905066
$ws.Range("M6").Value = 0.9639404333166532
$ws.Range("Q6").Value = 0.5296052411645183
$ws.Range("S6").Value = 0.52129602861432
$ws.Range("W6").Value = 0.8262346571285599

# Row 7 (M1) - boolean
$ws.Range("C7").Value = $true
$ws.Range("F7").Value = $true
$ws.Range("M7").Value = $true
$ws.Range("W7").Value = $true

# Row 8 (M3) - boolean
$ws.Range("C8").Value = $true
$ws.Range("F8").Value = $true
$ws.Range("I8").Value = $true
$ws.Range("M8").Value = $true
$ws.Range("S8").Value = $true
$ws.Range("W8").Value = $true

# Row 9 (M5) - boolean
$ws.Range("C9").Value = $true
$ws.Range("F9").Value = $true
$ws.Range("I9").Value = $true
$ws.Range("M9").Value = $true
$ws.Range("Q9").Value = $true
$ws.Range("S9").Value = $true
$ws.Range("W9").Value = $true

# Row 10 (position)
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("I10").Value = 2
$ws.Range("M10").Value = 1
$ws.Range("Q10").Value = 4
$ws.Range("S10").Value = 2
$ws.Range("W10").Value = 1
